# Add "snRNAseq-10xGenomics-v2" to the scrnaseq "assay_type list" sheet.
#
# The new value is inserted right after "scRNAseq-10xGenomics-v3" (row 3),
# pushing the existing "scRNAseq" / "sciRNAseq" / "snRNAseq" / "SNARE2-RNAseq"
# rows down by one. The data-validation list on the "Export as TSV" sheet
# (column L, assay_type) is widened to match the now-7-row list.

$wb = $excel.ActiveWorkbook

$assayTypeSheet = $wb.Worksheets.Item("assay_type list")

# Insert a new row 3 (shifting scRNAseq / sciRNAseq / snRNAseq / SNARE2-RNAseq
# down to rows 4-7) and populate it with the new assay type.
$assayTypeSheet.Rows("3:3").Insert()
$assayTypeSheet.Range("A3").Value = "snRNAseq-10xGenomics-v2"

# Widen the assay_type data validation on the main sheet from $A$1:$A$6 to
# $A$1:$A$7 so it covers the newly added row.
$mainSheet = $wb.Worksheets.Item("Export as TSV")
$assayTypeValidation = $mainSheet.Range("L2:L1048576").Validation
$assayTypeValidation.Formula1 = "='assay_type list'!`$A`$1:`$A`$7"
